$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Cells that simply get the "Y" mark (style stays the same, s="5")
$plainCells = @(
    "O4", "P4",
    "O5", "P5",
    "O6", "P6",
    "B7", "O7", "P7",
    "B8", "O8", "P8",
    "B9", "O9", "P9",
    "O10", "P10",
    "O11", "P11",
    "O12", "P12",
    "O13",
    "O14", "P14"
)

foreach ($addr in $plainCells) {
    $ws.Range($addr).Value = "Y"
}

# Cells that also pick up the highlighted fill (style goes from s="5" to s="8"),
# matching the formatting already used elsewhere in column B / row 13.
$ws.Range("B13").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("B11").PasteSpecial(-4122)
$ws.Range("B12").PasteSpecial(-4122)

$ws.Range("K13").Copy()
$ws.Range("P13").PasteSpecial(-4122)

$ws.Range("B10").Value = "Y"
$ws.Range("B11").Value = "Y"
$ws.Range("B12").Value = "Y"
$ws.Range("P13").Value = "Y"

# Update the last-known selection to match the author's final cursor position.
$ws.Range("Q14").Select()
